# Auto-generated edit script: updates horarios-141 schedule data
# for sheets LP1912, LP1912-215, 6203-6173 to the latest scrape.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- LP1912 ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 14:57:45"
$ws1.Cells.Item(3,1).Value = "Total filas: 252"
$ws1data = New-Object 'object[,]' 252,5
$ws1data[0,0] = "04:43:39"
$ws1data[0,1] = "04:45"
$ws1data[0,2] = "215A_EL PATO"
$ws1data[0,3] = 2
$ws1data[0,4] = "LP1912"
$ws1data[1,0] = "04:43:39"
$ws1data[1,1] = "04:53"
$ws1data[1,2] = "11_ETCHEVERRY"
$ws1data[1,3] = 10
$ws1data[1,4] = "LP1912"
$ws1data[2,0] = "04:56:49"
$ws1data[2,1] = "05:16"
$ws1data[2,2] = "17_ROMERO"
$ws1data[2,3] = 20
$ws1data[2,4] = "LP1912"
$ws1data[3,0] = "04:56:49"
$ws1data[3,1] = "05:22"
$ws1data[3,2] = "23_HERNANDEZ"
$ws1data[3,3] = 26
$ws1data[3,4] = "LP1912"
$ws1data[4,0] = "05:23:04"
$ws1data[4,1] = "05:23"
$ws1data[4,2] = "23_HERNANDEZ"
$ws1data[4,3] = 0
$ws1data[4,4] = "LP1912"
$ws1data[5,0] = "05:23:04"
$ws1data[5,1] = "05:32"
$ws1data[5,2] = "81_EL PELIGRO"
$ws1data[5,3] = 9
$ws1data[5,4] = "LP1912"
$ws1data[6,0] = "04:56:49"
$ws1data[6,1] = "05:34"
$ws1data[6,2] = "215B_EL PATO"
$ws1data[6,3] = 38
$ws1data[6,4] = "LP1912"
$ws1data[7,0] = "05:23:04"
$ws1data[7,1] = "05:44"
$ws1data[7,2] = "14_ABASTO"
$ws1data[7,3] = 21
$ws1data[7,4] = "LP1912"
$ws1data[8,0] = "04:56:49"
$ws1data[8,1] = "05:46"
$ws1data[8,2] = "15_ABASTO"
$ws1data[8,3] = 50
$ws1data[8,4] = "LP1912"
$ws1data[9,0] = "05:51:38"
$ws1data[9,1] = "05:52"
$ws1data[9,2] = "17_ROMERO"
$ws1data[9,3] = 1
$ws1data[9,4] = "LP1912"
$ws1data[10,0] = "04:56:49"
$ws1data[10,1] = "05:54"
$ws1data[10,2] = "10_OLMOS"
$ws1data[10,3] = 58
$ws1data[10,4] = "LP1912"
$ws1data[11,0] = "05:51:38"
$ws1data[11,1] = "06:03"
$ws1data[11,2] = "10_OLMOS"
$ws1data[11,3] = 12
$ws1data[11,4] = "LP1912"
$ws1data[12,0] = "04:56:49"
$ws1data[12,1] = "06:04"
$ws1data[12,2] = "16_SANTA ANA"
$ws1data[12,3] = 68
$ws1data[12,4] = "LP1912"
$ws1data[13,0] = "05:23:04"
$ws1data[13,1] = "06:04"
$ws1data[13,2] = "10_OLMOS"
$ws1data[13,3] = 41
$ws1data[13,4] = "LP1912"
$ws1data[14,0] = "05:51:38"
$ws1data[14,1] = "06:10"
$ws1data[14,2] = "215A_EL PATO"
$ws1data[14,3] = 19
$ws1data[14,4] = "LP1912"
$ws1data[15,0] = "05:23:04"
$ws1data[15,1] = "06:11"
$ws1data[15,2] = "215A_EL PATO"
$ws1data[15,3] = 48
$ws1data[15,4] = "LP1912"
$ws1data[16,0] = "04:56:49"
$ws1data[16,1] = "06:14"
$ws1data[16,2] = "225_HARAS DEL SUR"
$ws1data[16,3] = 78
$ws1data[16,4] = "LP1912"
$ws1data[17,0] = "04:56:49"
$ws1data[17,1] = "06:21"
$ws1data[17,2] = "26_HERNANDEZ"
$ws1data[17,3] = 85
$ws1data[17,4] = "LP1912"
$ws1data[18,0] = "06:19:29"
$ws1data[18,1] = "06:24"
$ws1data[18,2] = "11_ETCHEVERRY"
$ws1data[18,3] = 5
$ws1data[18,4] = "LP1912"
$ws1data[19,0] = "06:19:29"
$ws1data[19,1] = "06:27"
$ws1data[19,2] = "23_HERNANDEZ"
$ws1data[19,3] = 8
$ws1data[19,4] = "LP1912"
$ws1data[20,0] = "04:56:49"
$ws1data[20,1] = "06:29"
$ws1data[20,2] = "86_EST CHICA-ESC AGRARIA"
$ws1data[20,3] = 93
$ws1data[20,4] = "LP1912"
$ws1data[21,0] = "05:51:38"
$ws1data[21,1] = "06:30"
$ws1data[21,2] = "16_SANTA ANA"
$ws1data[21,3] = 39
$ws1data[21,4] = "LP1912"
$ws1data[22,0] = "06:19:29"
$ws1data[22,1] = "06:31"
$ws1data[22,2] = "17X38_ROMERO"
$ws1data[22,3] = 12
$ws1data[22,4] = "LP1912"
$ws1data[23,0] = "06:19:29"
$ws1data[23,1] = "06:31"
$ws1data[23,2] = "16_SANTA ANA"
$ws1data[23,3] = 12
$ws1data[23,4] = "LP1912"
$ws1data[24,0] = "06:19:29"
$ws1data[24,1] = "06:39"
$ws1data[24,2] = "225_C ROCA-H SUR"
$ws1data[24,3] = 20
$ws1data[24,4] = "LP1912"
$ws1data[25,0] = "04:56:49"
$ws1data[25,1] = "06:44"
$ws1data[25,2] = "225_C ROCA-H SUR"
$ws1data[25,3] = 108
$ws1data[25,4] = "LP1912"
$ws1data[26,0] = "04:56:49"
$ws1data[26,1] = "06:46"
$ws1data[26,2] = "215C_EL PATO"
$ws1data[26,3] = 110
$ws1data[26,4] = "LP1912"
$ws1data[27,0] = "05:51:38"
$ws1data[27,1] = "06:50"
$ws1data[27,2] = "215A_EL PATO"
$ws1data[27,3] = 59
$ws1data[27,4] = "LP1912"
$ws1data[28,0] = "06:46:06"
$ws1data[28,1] = "06:51"
$ws1data[28,2] = "215A_EL PATO"
$ws1data[28,3] = 5
$ws1data[28,4] = "LP1912"
$ws1data[29,0] = "06:46:06"
$ws1data[29,1] = "06:54"
$ws1data[29,2] = "14_ABASTO"
$ws1data[29,3] = 8
$ws1data[29,4] = "LP1912"
$ws1data[30,0] = "05:51:38"
$ws1data[30,1] = "07:00"
$ws1data[30,2] = "16_SANTA ANA"
$ws1data[30,3] = 69
$ws1data[30,4] = "LP1912"
$ws1data[31,0] = "06:58:01"
$ws1data[31,1] = "07:01"
$ws1data[31,2] = "16_SANTA ANA"
$ws1data[31,3] = 3
$ws1data[31,4] = "LP1912"
$ws1data[32,0] = "06:58:01"
$ws1data[32,1] = "07:04"
$ws1data[32,2] = "225_GOMEZ"
$ws1data[32,3] = 6
$ws1data[32,4] = "LP1912"
$ws1data[33,0] = "06:19:29"
$ws1data[33,1] = "07:06"
$ws1data[33,2] = "215C_EL PATO"
$ws1data[33,3] = 47
$ws1data[33,4] = "LP1912"
$ws1data[34,0] = "06:58:01"
$ws1data[34,1] = "07:07"
$ws1data[34,2] = "215C_EL PATO"
$ws1data[34,3] = 9
$ws1data[34,4] = "LP1912"
$ws1data[35,0] = "06:19:29"
$ws1data[35,1] = "07:13"
$ws1data[35,2] = "14X44_ABASTO"
$ws1data[35,3] = 54
$ws1data[35,4] = "LP1912"
$ws1data[36,0] = "06:58:01"
$ws1data[36,1] = "07:14"
$ws1data[36,2] = "14X44_ABASTO"
$ws1data[36,3] = 16
$ws1data[36,4] = "LP1912"
$ws1data[37,0] = "07:21:42"
$ws1data[37,1] = "07:21"
$ws1data[37,2] = "215A_EL PATO"
$ws1data[37,3] = 0
$ws1data[37,4] = "LP1912"
$ws1data[38,0] = "07:21:42"
$ws1data[38,1] = "07:23"
$ws1data[38,2] = "16_SANTA ANA"
$ws1data[38,3] = 2
$ws1data[38,4] = "LP1912"
$ws1data[39,0] = "06:58:01"
$ws1data[39,1] = "07:24"
$ws1data[39,2] = "16_SANTA ANA"
$ws1data[39,3] = 26
$ws1data[39,4] = "LP1912"
$ws1data[40,0] = "07:21:42"
$ws1data[40,1] = "07:29"
$ws1data[40,2] = "14_ABASTO"
$ws1data[40,3] = 8
$ws1data[40,4] = "LP1912"
$ws1data[41,0] = "07:21:42"
$ws1data[41,1] = "07:33"
$ws1data[41,2] = "23_HERNANDEZ"
$ws1data[41,3] = 12
$ws1data[41,4] = "LP1912"
$ws1data[42,0] = "06:58:01"
$ws1data[42,1] = "07:34"
$ws1data[42,2] = "23_HERNANDEZ"
$ws1data[42,3] = 36
$ws1data[42,4] = "LP1912"
$ws1data[43,0] = "07:21:42"
$ws1data[43,1] = "07:36"
$ws1data[43,2] = "17X38_ROMERO"
$ws1data[43,3] = 15
$ws1data[43,4] = "LP1912"
$ws1data[44,0] = "07:21:42"
$ws1data[44,1] = "07:36"
$ws1data[44,2] = "27_EL RETIRO"
$ws1data[44,3] = 15
$ws1data[44,4] = "LP1912"
$ws1data[45,0] = "06:58:01"
$ws1data[45,1] = "07:37"
$ws1data[45,2] = "27_EL RETIRO"
$ws1data[45,3] = 39
$ws1data[45,4] = "LP1912"
$ws1data[46,0] = "07:21:42"
$ws1data[46,1] = "07:41"
$ws1data[46,2] = "16_SANTA ANA"
$ws1data[46,3] = 20
$ws1data[46,4] = "LP1912"
$ws1data[47,0] = "07:21:42"
$ws1data[47,1] = "07:43"
$ws1data[47,2] = "10_OLMOS"
$ws1data[47,3] = 22
$ws1data[47,4] = "LP1912"
$ws1data[48,0] = "06:58:01"
$ws1data[48,1] = "07:44"
$ws1data[48,2] = "10_OLMOS"
$ws1data[48,3] = 46
$ws1data[48,4] = "LP1912"
$ws1data[49,0] = "07:48:14"
$ws1data[49,1] = "07:48"
$ws1data[49,2] = "215A_EL PATO"
$ws1data[49,3] = 0
$ws1data[49,4] = "LP1912"
$ws1data[50,0] = "07:21:42"
$ws1data[50,1] = "07:49"
$ws1data[50,2] = "15_ABASTO"
$ws1data[50,3] = 28
$ws1data[50,4] = "LP1912"
$ws1data[51,0] = "07:48:14"
$ws1data[51,1] = "07:58"
$ws1data[51,2] = "23_HERNANDEZ"
$ws1data[51,3] = 10
$ws1data[51,4] = "LP1912"
$ws1data[52,0] = "07:21:42"
$ws1data[52,1] = "07:59"
$ws1data[52,2] = "11_ETCHEVERRY"
$ws1data[52,3] = 38
$ws1data[52,4] = "LP1912"
$ws1data[53,0] = "08:00:50"
$ws1data[53,1] = "08:00"
$ws1data[53,2] = "215A_EL PATO"
$ws1data[53,3] = 0
$ws1data[53,4] = "LP1912"
$ws1data[54,0] = "08:00:50"
$ws1data[54,1] = "08:00"
$ws1data[54,2] = "11_ETCHEVERRY"
$ws1data[54,3] = 0
$ws1data[54,4] = "LP1912"
$ws1data[55,0] = "07:48:14"
$ws1data[55,1] = "08:01"
$ws1data[55,2] = "16_SANTA ANA"
$ws1data[55,3] = 13
$ws1data[55,4] = "LP1912"
$ws1data[56,0] = "06:46:06"
$ws1data[56,1] = "08:03"
$ws1data[56,2] = "23_HERNANDEZ"
$ws1data[56,3] = 77
$ws1data[56,4] = "LP1912"
$ws1data[57,0] = "08:00:50"
$ws1data[57,1] = "08:03"
$ws1data[57,2] = "17X38_ROMERO"
$ws1data[57,3] = 3
$ws1data[57,4] = "LP1912"
$ws1data[58,0] = "06:58:01"
$ws1data[58,1] = "08:04"
$ws1data[58,2] = "23_HERNANDEZ"
$ws1data[58,3] = 66
$ws1data[58,4] = "LP1912"
$ws1data[59,0] = "08:00:50"
$ws1data[59,1] = "08:14"
$ws1data[59,2] = "10_OLMOS"
$ws1data[59,3] = 14
$ws1data[59,4] = "LP1912"
$ws1data[60,0] = "08:00:50"
$ws1data[60,1] = "08:19"
$ws1data[60,2] = "15_ABASTO"
$ws1data[60,3] = 19
$ws1data[60,4] = "LP1912"
$ws1data[61,0] = "08:00:50"
$ws1data[61,1] = "08:21"
$ws1data[61,2] = "16_SANTA ANA"
$ws1data[61,3] = 21
$ws1data[61,4] = "LP1912"
$ws1data[62,0] = "08:00:50"
$ws1data[62,1] = "08:29"
$ws1data[62,2] = "14_ABASTO"
$ws1data[62,3] = 29
$ws1data[62,4] = "LP1912"
$ws1data[63,0] = "08:30:59"
$ws1data[63,1] = "08:30"
$ws1data[63,2] = "14_ABASTO"
$ws1data[63,3] = 0
$ws1data[63,4] = "LP1912"
$ws1data[64,0] = "08:30:59"
$ws1data[64,1] = "08:30"
$ws1data[64,2] = "215A_EL PATO"
$ws1data[64,3] = 0
$ws1data[64,4] = "LP1912"
$ws1data[65,0] = "08:30:59"
$ws1data[65,1] = "08:33"
$ws1data[65,2] = "215C_EL PATO"
$ws1data[65,3] = 3
$ws1data[65,4] = "LP1912"
$ws1data[66,0] = "08:00:50"
$ws1data[66,1] = "08:33"
$ws1data[66,2] = "23_HERNANDEZ"
$ws1data[66,3] = 33
$ws1data[66,4] = "LP1912"
$ws1data[67,0] = "07:48:14"
$ws1data[67,1] = "08:34"
$ws1data[67,2] = "215C_EL PATO"
$ws1data[67,3] = 46
$ws1data[67,4] = "LP1912"
$ws1data[68,0] = "08:30:59"
$ws1data[68,1] = "08:34"
$ws1data[68,2] = "23_HERNANDEZ"
$ws1data[68,3] = 4
$ws1data[68,4] = "LP1912"
$ws1data[69,0] = "08:30:59"
$ws1data[69,1] = "08:40"
$ws1data[69,2] = "16_SANTA ANA"
$ws1data[69,3] = 10
$ws1data[69,4] = "LP1912"
$ws1data[70,0] = "08:00:50"
$ws1data[70,1] = "08:41"
$ws1data[70,2] = "16_SANTA ANA"
$ws1data[70,3] = 41
$ws1data[70,4] = "LP1912"
$ws1data[71,0] = "08:00:50"
$ws1data[71,1] = "08:44"
$ws1data[71,2] = "10_OLMOS"
$ws1data[71,3] = 44
$ws1data[71,4] = "LP1912"
$ws1data[72,0] = "08:30:59"
$ws1data[72,1] = "08:47"
$ws1data[72,2] = "215A_EL PATO"
$ws1data[72,3] = 17
$ws1data[72,4] = "LP1912"
$ws1data[73,0] = "08:48:29"
$ws1data[73,1] = "08:48"
$ws1data[73,2] = "215A_EL PATO"
$ws1data[73,3] = 0
$ws1data[73,4] = "LP1912"
$ws1data[74,0] = "08:30:59"
$ws1data[74,1] = "08:50"
$ws1data[74,2] = "10_OLMOS"
$ws1data[74,3] = 20
$ws1data[74,4] = "LP1912"
$ws1data[75,0] = "07:21:42"
$ws1data[75,1] = "08:51"
$ws1data[75,2] = "16_P MOR-SANTA ANA"
$ws1data[75,3] = 90
$ws1data[75,4] = "LP1912"
$ws1data[76,0] = "08:48:29"
$ws1data[76,1] = "08:51"
$ws1data[76,2] = "10_OLMOS"
$ws1data[76,3] = 3
$ws1data[76,4] = "LP1912"
$ws1data[77,0] = "07:48:14"
$ws1data[77,1] = "08:52"
$ws1data[77,2] = "16_P MOR-SANTA ANA"
$ws1data[77,3] = 64
$ws1data[77,4] = "LP1912"
$ws1data[78,0] = "08:00:50"
$ws1data[78,1] = "08:55"
$ws1data[78,2] = "16_P MOR-SANTA ANA"
$ws1data[78,3] = 55
$ws1data[78,4] = "LP1912"
$ws1data[79,0] = "08:56:14"
$ws1data[79,1] = "08:56"
$ws1data[79,2] = "215B_EL PATO"
$ws1data[79,3] = 0
$ws1data[79,4] = "LP1912"
$ws1data[80,0] = "08:56:14"
$ws1data[80,1] = "08:56"
$ws1data[80,2] = "215A_EL PATO"
$ws1data[80,3] = 0
$ws1data[80,4] = "LP1912"
$ws1data[81,0] = "08:48:29"
$ws1data[81,1] = "08:59"
$ws1data[81,2] = "215B_EL PATO"
$ws1data[81,3] = 11
$ws1data[81,4] = "LP1912"
$ws1data[82,0] = "08:30:59"
$ws1data[82,1] = "08:59"
$ws1data[82,2] = "16_P MOR-SANTA ANA"
$ws1data[82,3] = 29
$ws1data[82,4] = "LP1912"
$ws1data[83,0] = "07:48:14"
$ws1data[83,1] = "09:00"
$ws1data[83,2] = "215B_EL PATO"
$ws1data[83,3] = 72
$ws1data[83,4] = "LP1912"
$ws1data[84,0] = "08:30:59"
$ws1data[84,1] = "09:00"
$ws1data[84,2] = "16_SANTA ANA"
$ws1data[84,3] = 30
$ws1data[84,4] = "LP1912"
$ws1data[85,0] = "08:48:29"
$ws1data[85,1] = "09:01"
$ws1data[85,2] = "16_P MOR-SANTA ANA"
$ws1data[85,3] = 13
$ws1data[85,4] = "LP1912"
$ws1data[86,0] = "08:56:14"
$ws1data[86,1] = "09:01"
$ws1data[86,2] = "16_SANTA ANA"
$ws1data[86,3] = 5
$ws1data[86,4] = "LP1912"
$ws1data[87,0] = "08:30:59"
$ws1data[87,1] = "09:02"
$ws1data[87,2] = "17X38_ROMERO"
$ws1data[87,3] = 32
$ws1data[87,4] = "LP1912"
$ws1data[88,0] = "08:56:14"
$ws1data[88,1] = "09:03"
$ws1data[88,2] = "17X38_ROMERO"
$ws1data[88,3] = 7
$ws1data[88,4] = "LP1912"
$ws1data[89,0] = "08:56:14"
$ws1data[89,1] = "09:03"
$ws1data[89,2] = "23_HERNANDEZ"
$ws1data[89,3] = 7
$ws1data[89,4] = "LP1912"
$ws1data[90,0] = "08:48:29"
$ws1data[90,1] = "09:14"
$ws1data[90,2] = "11_ETCHEVERRY"
$ws1data[90,3] = 26
$ws1data[90,4] = "LP1912"
$ws1data[91,0] = "08:00:50"
$ws1data[91,1] = "09:14"
$ws1data[91,2] = "27_EL RETIRO"
$ws1data[91,3] = 74
$ws1data[91,4] = "LP1912"
$ws1data[92,0] = "08:56:14"
$ws1data[92,1] = "09:15"
$ws1data[92,2] = "11_ETCHEVERRY"
$ws1data[92,3] = 19
$ws1data[92,4] = "LP1912"
$ws1data[93,0] = "08:48:29"
$ws1data[93,1] = "09:16"
$ws1data[93,2] = "27_EL RETIRO"
$ws1data[93,3] = 28
$ws1data[93,4] = "LP1912"
$ws1data[94,0] = "08:56:14"
$ws1data[94,1] = "09:17"
$ws1data[94,2] = "27_EL RETIRO"
$ws1data[94,3] = 21
$ws1data[94,4] = "LP1912"
$ws1data[95,0] = "08:30:59"
$ws1data[95,1] = "09:18"
$ws1data[95,2] = "215_EL PELIGRO"
$ws1data[95,3] = 48
$ws1data[95,4] = "LP1912"
$ws1data[96,0] = "08:56:14"
$ws1data[96,1] = "09:19"
$ws1data[96,2] = "215_EL PELIGRO"
$ws1data[96,3] = 23
$ws1data[96,4] = "LP1912"
$ws1data[97,0] = "08:56:14"
$ws1data[97,1] = "09:21"
$ws1data[97,2] = "16_SANTA ANA"
$ws1data[97,3] = 25
$ws1data[97,4] = "LP1912"
$ws1data[98,0] = "08:56:14"
$ws1data[98,1] = "09:29"
$ws1data[98,2] = "10_OLMOS"
$ws1data[98,3] = 33
$ws1data[98,4] = "LP1912"
$ws1data[99,0] = "09:31:15"
$ws1data[99,1] = "09:31"
$ws1data[99,2] = "23_HERNANDEZ"
$ws1data[99,3] = 0
$ws1data[99,4] = "LP1912"
$ws1data[100,0] = "08:56:14"
$ws1data[100,1] = "09:33"
$ws1data[100,2] = "23_HERNANDEZ"
$ws1data[100,3] = 37
$ws1data[100,4] = "LP1912"
$ws1data[101,0] = "08:30:59"
$ws1data[101,1] = "09:33"
$ws1data[101,2] = "15_ABASTO"
$ws1data[101,3] = 63
$ws1data[101,4] = "LP1912"
$ws1data[102,0] = "09:31:15"
$ws1data[102,1] = "09:34"
$ws1data[102,2] = "15_ABASTO"
$ws1data[102,3] = 3
$ws1data[102,4] = "LP1912"
$ws1data[103,0] = "09:31:15"
$ws1data[103,1] = "09:41"
$ws1data[103,2] = "16_SANTA ANA"
$ws1data[103,3] = 10
$ws1data[103,4] = "LP1912"
$ws1data[104,0] = "08:48:29"
$ws1data[104,1] = "09:44"
$ws1data[104,2] = "14_ABASTO"
$ws1data[104,3] = 56
$ws1data[104,4] = "LP1912"
$ws1data[105,0] = "09:31:15"
$ws1data[105,1] = "09:45"
$ws1data[105,2] = "14_ABASTO"
$ws1data[105,3] = 14
$ws1data[105,4] = "LP1912"
$ws1data[106,0] = "08:30:59"
$ws1data[106,1] = "09:48"
$ws1data[106,2] = "15_ABASTO"
$ws1data[106,3] = 78
$ws1data[106,4] = "LP1912"
$ws1data[107,0] = "08:30:59"
$ws1data[107,1] = "09:50"
$ws1data[107,2] = "16_P MOR-SANTA ANA"
$ws1data[107,3] = 80
$ws1data[107,4] = "LP1912"
$ws1data[108,0] = "09:31:15"
$ws1data[108,1] = "09:51"
$ws1data[108,2] = "16_P MOR-SANTA ANA"
$ws1data[108,3] = 20
$ws1data[108,4] = "LP1912"
$ws1data[109,0] = "08:30:59"
$ws1data[109,1] = "09:55"
$ws1data[109,2] = "10_OLMOS"
$ws1data[109,3] = 85
$ws1data[109,4] = "LP1912"
$ws1data[110,0] = "09:31:15"
$ws1data[110,1] = "09:56"
$ws1data[110,2] = "10_OLMOS"
$ws1data[110,3] = 25
$ws1data[110,4] = "LP1912"
$ws1data[111,0] = "09:31:15"
$ws1data[111,1] = "10:01"
$ws1data[111,2] = "16_SANTA ANA"
$ws1data[111,3] = 30
$ws1data[111,4] = "LP1912"
$ws1data[112,0] = "08:30:59"
$ws1data[112,1] = "10:03"
$ws1data[112,2] = "215C_EL PATO"
$ws1data[112,3] = 93
$ws1data[112,4] = "LP1912"
$ws1data[113,0] = "09:31:15"
$ws1data[113,1] = "10:03"
$ws1data[113,2] = "23_HERNANDEZ"
$ws1data[113,3] = 32
$ws1data[113,4] = "LP1912"
$ws1data[114,0] = "09:31:15"
$ws1data[114,1] = "10:04"
$ws1data[114,2] = "215C_EL PATO"
$ws1data[114,3] = 33
$ws1data[114,4] = "LP1912"
$ws1data[115,0] = "09:31:15"
$ws1data[115,1] = "10:08"
$ws1data[115,2] = "11_ETCHEVERRY"
$ws1data[115,3] = 37
$ws1data[115,4] = "LP1912"
$ws1data[116,0] = "09:31:15"
$ws1data[116,1] = "10:19"
$ws1data[116,2] = "17_ROMERO"
$ws1data[116,3] = 48
$ws1data[116,4] = "LP1912"
$ws1data[117,0] = "08:30:59"
$ws1data[117,1] = "10:19"
$ws1data[117,2] = "10_OLMOS"
$ws1data[117,3] = 109
$ws1data[117,4] = "LP1912"
$ws1data[118,0] = "09:31:15"
$ws1data[118,1] = "10:20"
$ws1data[118,2] = "10_OLMOS"
$ws1data[118,3] = 49
$ws1data[118,4] = "LP1912"
$ws1data[119,0] = "10:26:25"
$ws1data[119,1] = "10:32"
$ws1data[119,2] = "14_ABASTO"
$ws1data[119,3] = 6
$ws1data[119,4] = "LP1912"
$ws1data[120,0] = "08:56:14"
$ws1data[120,1] = "10:33"
$ws1data[120,2] = "14_ABASTO"
$ws1data[120,3] = 97
$ws1data[120,4] = "LP1912"
$ws1data[121,0] = "10:26:25"
$ws1data[121,1] = "10:33"
$ws1data[121,2] = "23_HERNANDEZ"
$ws1data[121,3] = 7
$ws1data[121,4] = "LP1912"
$ws1data[122,0] = "09:31:15"
$ws1data[122,1] = "10:34"
$ws1data[122,2] = "14_ABASTO"
$ws1data[122,3] = 63
$ws1data[122,4] = "LP1912"
$ws1data[123,0] = "10:26:25"
$ws1data[123,1] = "10:34"
$ws1data[123,2] = "15_ABASTO"
$ws1data[123,3] = 8
$ws1data[123,4] = "LP1912"
$ws1data[124,0] = "08:48:29"
$ws1data[124,1] = "10:36"
$ws1data[124,2] = "14_ABASTO"
$ws1data[124,3] = 108
$ws1data[124,4] = "LP1912"
$ws1data[125,0] = "10:26:25"
$ws1data[125,1] = "10:41"
$ws1data[125,2] = "16_SANTA ANA"
$ws1data[125,3] = 15
$ws1data[125,4] = "LP1912"
$ws1data[126,0] = "10:26:25"
$ws1data[126,1] = "10:44"
$ws1data[126,2] = "10_OLMOS"
$ws1data[126,3] = 18
$ws1data[126,4] = "LP1912"
$ws1data[127,0] = "10:26:25"
$ws1data[127,1] = "10:49"
$ws1data[127,2] = "15_ABASTO"
$ws1data[127,3] = 23
$ws1data[127,4] = "LP1912"
$ws1data[128,0] = "10:26:25"
$ws1data[128,1] = "10:51"
$ws1data[128,2] = "16_P MOR-SANTA ANA"
$ws1data[128,3] = 25
$ws1data[128,4] = "LP1912"
$ws1data[129,0] = "10:55:35"
$ws1data[129,1] = "10:56"
$ws1data[129,2] = "14_ABASTO"
$ws1data[129,3] = 1
$ws1data[129,4] = "LP1912"
$ws1data[130,0] = "10:55:35"
$ws1data[130,1] = "10:56"
$ws1data[130,2] = "27_EL RETIRO"
$ws1data[130,3] = 1
$ws1data[130,4] = "LP1912"
$ws1data[131,0] = "10:26:25"
$ws1data[131,1] = "10:57"
$ws1data[131,2] = "27_EL RETIRO"
$ws1data[131,3] = 31
$ws1data[131,4] = "LP1912"
$ws1data[132,0] = "10:55:35"
$ws1data[132,1] = "11:01"
$ws1data[132,2] = "16_SANTA ANA"
$ws1data[132,3] = 6
$ws1data[132,4] = "LP1912"
$ws1data[133,0] = "10:55:35"
$ws1data[133,1] = "11:03"
$ws1data[133,2] = "23_HERNANDEZ"
$ws1data[133,3] = 8
$ws1data[133,4] = "LP1912"
$ws1data[134,0] = "10:55:35"
$ws1data[134,1] = "11:04"
$ws1data[134,2] = "17_ROMERO"
$ws1data[134,3] = 9
$ws1data[134,4] = "LP1912"
$ws1data[135,0] = "10:55:35"
$ws1data[135,1] = "11:08"
$ws1data[135,2] = "225_C ROCA-H SUR"
$ws1data[135,3] = 13
$ws1data[135,4] = "LP1912"
$ws1data[136,0] = "11:17:16"
$ws1data[136,1] = "11:18"
$ws1data[136,2] = "17_ROMERO"
$ws1data[136,3] = 1
$ws1data[136,4] = "LP1912"
$ws1data[137,0] = "10:55:35"
$ws1data[137,1] = "11:19"
$ws1data[137,2] = "215C_EL PATO"
$ws1data[137,3] = 24
$ws1data[137,4] = "LP1912"
$ws1data[138,0] = "10:55:35"
$ws1data[138,1] = "11:20"
$ws1data[138,2] = "11_ETCHEVERRY"
$ws1data[138,3] = 25
$ws1data[138,4] = "LP1912"
$ws1data[139,0] = "11:17:16"
$ws1data[139,1] = "11:21"
$ws1data[139,2] = "11_ETCHEVERRY"
$ws1data[139,3] = 4
$ws1data[139,4] = "LP1912"
$ws1data[140,0] = "11:17:16"
$ws1data[140,1] = "11:33"
$ws1data[140,2] = "215A_EL PATO"
$ws1data[140,3] = 16
$ws1data[140,4] = "LP1912"
$ws1data[141,0] = "11:17:16"
$ws1data[141,1] = "11:33"
$ws1data[141,2] = "23_HERNANDEZ"
$ws1data[141,3] = 16
$ws1data[141,4] = "LP1912"
$ws1data[142,0] = "11:17:16"
$ws1data[142,1] = "11:41"
$ws1data[142,2] = "16_SANTA ANA"
$ws1data[142,3] = 24
$ws1data[142,4] = "LP1912"
$ws1data[143,0] = "11:44:58"
$ws1data[143,1] = "11:44"
$ws1data[143,2] = "215B_EL PATO"
$ws1data[143,3] = 0
$ws1data[143,4] = "LP1912"
$ws1data[144,0] = "11:17:16"
$ws1data[144,1] = "11:45"
$ws1data[144,2] = "215B_EL PATO"
$ws1data[144,3] = 28
$ws1data[144,4] = "LP1912"
$ws1data[145,0] = "11:44:58"
$ws1data[145,1] = "11:45"
$ws1data[145,2] = "17_ROMERO"
$ws1data[145,3] = 1
$ws1data[145,4] = "LP1912"
$ws1data[146,0] = "11:44:58"
$ws1data[146,1] = "11:48"
$ws1data[146,2] = "15_ABASTO"
$ws1data[146,3] = 4
$ws1data[146,4] = "LP1912"
$ws1data[147,0] = "11:17:16"
$ws1data[147,1] = "11:49"
$ws1data[147,2] = "15_ABASTO"
$ws1data[147,3] = 32
$ws1data[147,4] = "LP1912"
$ws1data[148,0] = "11:44:58"
$ws1data[148,1] = "11:50"
$ws1data[148,2] = "16_P MOR-SANTA ANA"
$ws1data[148,3] = 6
$ws1data[148,4] = "LP1912"
$ws1data[149,0] = "11:17:16"
$ws1data[149,1] = "11:51"
$ws1data[149,2] = "16_P MOR-SANTA ANA"
$ws1data[149,3] = 34
$ws1data[149,4] = "LP1912"
$ws1data[150,0] = "11:44:58"
$ws1data[150,1] = "11:55"
$ws1data[150,2] = "225_GOMEZ"
$ws1data[150,3] = 11
$ws1data[150,4] = "LP1912"
$ws1data[151,0] = "11:17:16"
$ws1data[151,1] = "11:56"
$ws1data[151,2] = "225_GOMEZ"
$ws1data[151,3] = 39
$ws1data[151,4] = "LP1912"
$ws1data[152,0] = "11:58:47"
$ws1data[152,1] = "11:58"
$ws1data[152,2] = "16_SANTA ANA"
$ws1data[152,3] = 0
$ws1data[152,4] = "LP1912"
$ws1data[153,0] = "11:58:47"
$ws1data[153,1] = "12:00"
$ws1data[153,2] = "17_ROMERO"
$ws1data[153,3] = 2
$ws1data[153,4] = "LP1912"
$ws1data[154,0] = "11:17:16"
$ws1data[154,1] = "12:01"
$ws1data[154,2] = "23_HERNANDEZ"
$ws1data[154,3] = 44
$ws1data[154,4] = "LP1912"
$ws1data[155,0] = "11:44:58"
$ws1data[155,1] = "12:01"
$ws1data[155,2] = "16_SANTA ANA"
$ws1data[155,3] = 17
$ws1data[155,4] = "LP1912"
$ws1data[156,0] = "11:44:58"
$ws1data[156,1] = "12:03"
$ws1data[156,2] = "17_ROMERO"
$ws1data[156,3] = 19
$ws1data[156,4] = "LP1912"
$ws1data[157,0] = "11:58:47"
$ws1data[157,1] = "12:03"
$ws1data[157,2] = "23_HERNANDEZ"
$ws1data[157,3] = 5
$ws1data[157,4] = "LP1912"
$ws1data[158,0] = "11:58:47"
$ws1data[158,1] = "12:04"
$ws1data[158,2] = "17_ROMERO"
$ws1data[158,3] = 6
$ws1data[158,4] = "LP1912"
$ws1data[159,0] = "11:44:58"
$ws1data[159,1] = "12:04"
$ws1data[159,2] = "23_HERNANDEZ"
$ws1data[159,3] = 20
$ws1data[159,4] = "LP1912"
$ws1data[160,0] = "11:44:58"
$ws1data[160,1] = "12:08"
$ws1data[160,2] = "14_ABASTO"
$ws1data[160,3] = 24
$ws1data[160,4] = "LP1912"
$ws1data[161,0] = "11:17:16"
$ws1data[161,1] = "12:09"
$ws1data[161,2] = "14_ABASTO"
$ws1data[161,3] = 52
$ws1data[161,4] = "LP1912"
$ws1data[162,0] = "11:44:58"
$ws1data[162,1] = "12:18"
$ws1data[162,2] = "15_ABASTO"
$ws1data[162,3] = 34
$ws1data[162,4] = "LP1912"
$ws1data[163,0] = "11:44:58"
$ws1data[163,1] = "12:19"
$ws1data[163,2] = "10_OLMOS"
$ws1data[163,3] = 35
$ws1data[163,4] = "LP1912"
$ws1data[164,0] = "11:58:47"
$ws1data[164,1] = "12:19"
$ws1data[164,2] = "15_ABASTO"
$ws1data[164,3] = 21
$ws1data[164,4] = "LP1912"
$ws1data[165,0] = "11:58:47"
$ws1data[165,1] = "12:20"
$ws1data[165,2] = "10_OLMOS"
$ws1data[165,3] = 22
$ws1data[165,4] = "LP1912"
$ws1data[166,0] = "11:44:58"
$ws1data[166,1] = "12:21"
$ws1data[166,2] = "16_SANTA ANA"
$ws1data[166,3] = 37
$ws1data[166,4] = "LP1912"
$ws1data[167,0] = "12:23:39"
$ws1data[167,1] = "12:23"
$ws1data[167,2] = "17_ROMERO"
$ws1data[167,3] = 0
$ws1data[167,4] = "LP1912"
$ws1data[168,0] = "12:23:39"
$ws1data[168,1] = "12:29"
$ws1data[168,2] = "10_OLMOS"
$ws1data[168,3] = 6
$ws1data[168,4] = "LP1912"
$ws1data[169,0] = "12:23:39"
$ws1data[169,1] = "12:32"
$ws1data[169,2] = "11_ETCHEVERRY"
$ws1data[169,3] = 9
$ws1data[169,4] = "LP1912"
$ws1data[170,0] = "11:17:16"
$ws1data[170,1] = "12:33"
$ws1data[170,2] = "11_ETCHEVERRY"
$ws1data[170,3] = 76
$ws1data[170,4] = "LP1912"
$ws1data[171,0] = "12:23:39"
$ws1data[171,1] = "12:33"
$ws1data[171,2] = "215C_EL PATO"
$ws1data[171,3] = 10
$ws1data[171,4] = "LP1912"
$ws1data[172,0] = "12:23:39"
$ws1data[172,1] = "12:33"
$ws1data[172,2] = "23_HERNANDEZ"
$ws1data[172,3] = 10
$ws1data[172,4] = "LP1912"
$ws1data[173,0] = "11:17:16"
$ws1data[173,1] = "12:34"
$ws1data[173,2] = "215C_EL PATO"
$ws1data[173,3] = 77
$ws1data[173,4] = "LP1912"
$ws1data[174,0] = "12:23:39"
$ws1data[174,1] = "12:36"
$ws1data[174,2] = "27_EL RETIRO"
$ws1data[174,3] = 13
$ws1data[174,4] = "LP1912"
$ws1data[175,0] = "12:43:41"
$ws1data[175,1] = "12:43"
$ws1data[175,2] = "16_SANTA ANA"
$ws1data[175,3] = 0
$ws1data[175,4] = "LP1912"
$ws1data[176,0] = "12:43:41"
$ws1data[176,1] = "12:44"
$ws1data[176,2] = "17_ROMERO"
$ws1data[176,3] = 1
$ws1data[176,4] = "LP1912"
$ws1data[177,0] = "11:44:58"
$ws1data[177,1] = "12:46"
$ws1data[177,2] = "10_OLMOS"
$ws1data[177,3] = 62
$ws1data[177,4] = "LP1912"
$ws1data[178,0] = "12:43:41"
$ws1data[178,1] = "12:47"
$ws1data[178,2] = "10_OLMOS"
$ws1data[178,3] = 4
$ws1data[178,4] = "LP1912"
$ws1data[179,0] = "12:43:41"
$ws1data[179,1] = "12:49"
$ws1data[179,2] = "15_ABASTO"
$ws1data[179,3] = 6
$ws1data[179,4] = "LP1912"
$ws1data[180,0] = "11:44:58"
$ws1data[180,1] = "12:50"
$ws1data[180,2] = "16_P MOR-SANTA ANA"
$ws1data[180,3] = 66
$ws1data[180,4] = "LP1912"
$ws1data[181,0] = "12:43:41"
$ws1data[181,1] = "12:51"
$ws1data[181,2] = "16_P MOR-SANTA ANA"
$ws1data[181,3] = 8
$ws1data[181,4] = "LP1912"
$ws1data[182,0] = "12:57:23"
$ws1data[182,1] = "12:58"
$ws1data[182,2] = "17_ROMERO"
$ws1data[182,3] = 1
$ws1data[182,4] = "LP1912"
$ws1data[183,0] = "11:17:16"
$ws1data[183,1] = "13:00"
$ws1data[183,2] = "14_ABASTO"
$ws1data[183,3] = 103
$ws1data[183,4] = "LP1912"
$ws1data[184,0] = "12:57:23"
$ws1data[184,1] = "13:02"
$ws1data[184,2] = "16_P MOR-SANTA ANA"
$ws1data[184,3] = 5
$ws1data[184,4] = "LP1912"
$ws1data[185,0] = "12:57:23"
$ws1data[185,1] = "13:02"
$ws1data[185,2] = "14_ABASTO"
$ws1data[185,3] = 5
$ws1data[185,4] = "LP1912"
$ws1data[186,0] = "12:57:23"
$ws1data[186,1] = "13:03"
$ws1data[186,2] = "23_HERNANDEZ"
$ws1data[186,3] = 6
$ws1data[186,4] = "LP1912"
$ws1data[187,0] = "12:43:41"
$ws1data[187,1] = "13:03"
$ws1data[187,2] = "215C_EL PATO"
$ws1data[187,3] = 20
$ws1data[187,4] = "LP1912"
$ws1data[188,0] = "11:58:47"
$ws1data[188,1] = "13:04"
$ws1data[188,2] = "23_HERNANDEZ"
$ws1data[188,3] = 66
$ws1data[188,4] = "LP1912"
$ws1data[189,0] = "12:57:23"
$ws1data[189,1] = "13:04"
$ws1data[189,2] = "215C_EL PATO"
$ws1data[189,3] = 7
$ws1data[189,4] = "LP1912"
$ws1data[190,0] = "11:44:58"
$ws1data[190,1] = "13:16"
$ws1data[190,2] = "11X44_P_ETCHEVERRY"
$ws1data[190,3] = 92
$ws1data[190,4] = "LP1912"
$ws1data[191,0] = "12:57:23"
$ws1data[191,1] = "13:17"
$ws1data[191,2] = "11X44_P_ETCHEVERRY"
$ws1data[191,3] = 20
$ws1data[191,4] = "LP1912"
$ws1data[192,0] = "12:43:41"
$ws1data[192,1] = "13:18"
$ws1data[192,2] = "215_ALUAR"
$ws1data[192,3] = 35
$ws1data[192,4] = "LP1912"
$ws1data[193,0] = "12:57:23"
$ws1data[193,1] = "13:19"
$ws1data[193,2] = "215_ALUAR"
$ws1data[193,3] = 22
$ws1data[193,4] = "LP1912"
$ws1data[194,0] = "12:57:23"
$ws1data[194,1] = "13:21"
$ws1data[194,2] = "16_SANTA ANA"
$ws1data[194,3] = 24
$ws1data[194,4] = "LP1912"
$ws1data[195,0] = "13:22:41"
$ws1data[195,1] = "13:27"
$ws1data[195,2] = "16_P MOR-SANTA ANA"
$ws1data[195,3] = 5
$ws1data[195,4] = "LP1912"
$ws1data[196,0] = "13:22:41"
$ws1data[196,1] = "13:32"
$ws1data[196,2] = "14_ABASTO"
$ws1data[196,3] = 10
$ws1data[196,4] = "LP1912"
$ws1data[197,0] = "13:22:41"
$ws1data[197,1] = "13:33"
$ws1data[197,2] = "23_HERNANDEZ"
$ws1data[197,3] = 11
$ws1data[197,4] = "LP1912"
$ws1data[198,0] = "11:58:47"
$ws1data[198,1] = "13:37"
$ws1data[198,2] = "215A_EL PATO"
$ws1data[198,3] = 99
$ws1data[198,4] = "LP1912"
$ws1data[199,0] = "13:22:41"
$ws1data[199,1] = "13:38"
$ws1data[199,2] = "215A_EL PATO"
$ws1data[199,3] = 16
$ws1data[199,4] = "LP1912"
$ws1data[200,0] = "13:22:41"
$ws1data[200,1] = "13:41"
$ws1data[200,2] = "16_SANTA ANA"
$ws1data[200,3] = 19
$ws1data[200,4] = "LP1912"
$ws1data[201,0] = "13:22:41"
$ws1data[201,1] = "13:47"
$ws1data[201,2] = "225_GOMEZ"
$ws1data[201,3] = 25
$ws1data[201,4] = "LP1912"
$ws1data[202,0] = "13:51:40"
$ws1data[202,1] = "13:51"
$ws1data[202,2] = "16_P MOR-SANTA ANA"
$ws1data[202,3] = 0
$ws1data[202,4] = "LP1912"
$ws1data[203,0] = "13:51:40"
$ws1data[203,1] = "13:59"
$ws1data[203,2] = "17_ROMERO"
$ws1data[203,3] = 8
$ws1data[203,4] = "LP1912"
$ws1data[204,0] = "13:51:40"
$ws1data[204,1] = "14:02"
$ws1data[204,2] = "14_ABASTO"
$ws1data[204,3] = 11
$ws1data[204,4] = "LP1912"
$ws1data[205,0] = "13:51:40"
$ws1data[205,1] = "14:03"
$ws1data[205,2] = "23_HERNANDEZ"
$ws1data[205,3] = 12
$ws1data[205,4] = "LP1912"
$ws1data[206,0] = "13:22:41"
$ws1data[206,1] = "14:13"
$ws1data[206,2] = "23_HERNANDEZ"
$ws1data[206,3] = 51
$ws1data[206,4] = "LP1912"
$ws1data[207,0] = "14:15:51"
$ws1data[207,1] = "14:16"
$ws1data[207,2] = "27_EL RETIRO"
$ws1data[207,3] = 1
$ws1data[207,4] = "LP1912"
$ws1data[208,0] = "14:15:51"
$ws1data[208,1] = "14:16"
$ws1data[208,2] = "10_OLMOS"
$ws1data[208,3] = 1
$ws1data[208,4] = "LP1912"
$ws1data[209,0] = "12:57:23"
$ws1data[209,1] = "14:17"
$ws1data[209,2] = "27_EL RETIRO"
$ws1data[209,3] = 80
$ws1data[209,4] = "LP1912"
$ws1data[210,0] = "13:51:40"
$ws1data[210,1] = "14:17"
$ws1data[210,2] = "10_OLMOS"
$ws1data[210,3] = 26
$ws1data[210,4] = "LP1912"
$ws1data[211,0] = "14:15:51"
$ws1data[211,1] = "14:18"
$ws1data[211,2] = "215C_EL PATO"
$ws1data[211,3] = 3
$ws1data[211,4] = "LP1912"
$ws1data[212,0] = "12:57:23"
$ws1data[212,1] = "14:19"
$ws1data[212,2] = "215C_EL PATO"
$ws1data[212,3] = 82
$ws1data[212,4] = "LP1912"
$ws1data[213,0] = "13:51:40"
$ws1data[213,1] = "14:20"
$ws1data[213,2] = "16_SANTA ANA"
$ws1data[213,3] = 29
$ws1data[213,4] = "LP1912"
$ws1data[214,0] = "14:15:51"
$ws1data[214,1] = "14:21"
$ws1data[214,2] = "16_SANTA ANA"
$ws1data[214,3] = 6
$ws1data[214,4] = "LP1912"
$ws1data[215,0] = "14:15:51"
$ws1data[215,1] = "14:32"
$ws1data[215,2] = "11_ETCHEVERRY"
$ws1data[215,3] = 17
$ws1data[215,4] = "LP1912"
$ws1data[216,0] = "14:15:51"
$ws1data[216,1] = "14:33"
$ws1data[216,2] = "23_HERNANDEZ"
$ws1data[216,3] = 18
$ws1data[216,4] = "LP1912"
$ws1data[217,0] = "14:15:51"
$ws1data[217,1] = "14:38"
$ws1data[217,2] = "15X38_ABASTO"
$ws1data[217,3] = 23
$ws1data[217,4] = "LP1912"
$ws1data[218,0] = "13:51:40"
$ws1data[218,1] = "14:39"
$ws1data[218,2] = "23_HERNANDEZ"
$ws1data[218,3] = 48
$ws1data[218,4] = "LP1912"
$ws1data[219,0] = "14:42:53"
$ws1data[219,1] = "14:42"
$ws1data[219,2] = "16_SANTA ANA"
$ws1data[219,3] = 0
$ws1data[219,4] = "LP1912"
$ws1data[220,0] = "14:42:53"
$ws1data[220,1] = "14:46"
$ws1data[220,2] = "10_OLMOS"
$ws1data[220,3] = 4
$ws1data[220,4] = "LP1912"
$ws1data[221,0] = "13:51:40"
$ws1data[221,1] = "14:47"
$ws1data[221,2] = "10_OLMOS"
$ws1data[221,3] = 56
$ws1data[221,4] = "LP1912"
$ws1data[222,0] = "14:42:53"
$ws1data[222,1] = "14:51"
$ws1data[222,2] = "16_P MOR-SANTA ANA"
$ws1data[222,3] = 9
$ws1data[222,4] = "LP1912"
$ws1data[223,0] = "14:57:45"
$ws1data[223,1] = "14:57"
$ws1data[223,2] = "16_SANTA ANA"
$ws1data[223,3] = 0
$ws1data[223,4] = "LP1912"
$ws1data[224,0] = "14:57:45"
$ws1data[224,1] = "14:58"
$ws1data[224,2] = "215B_EL PATO"
$ws1data[224,3] = 1
$ws1data[224,4] = "LP1912"
$ws1data[225,0] = "14:42:53"
$ws1data[225,1] = "14:59"
$ws1data[225,2] = "215B_EL PATO"
$ws1data[225,3] = 17
$ws1data[225,4] = "LP1912"
$ws1data[226,0] = "14:57:45"
$ws1data[226,1] = "14:59"
$ws1data[226,2] = "81_EL PELIGRO"
$ws1data[226,3] = 2
$ws1data[226,4] = "LP1912"
$ws1data[227,0] = "14:57:45"
$ws1data[227,1] = "15:03"
$ws1data[227,2] = "23_HERNANDEZ"
$ws1data[227,3] = 6
$ws1data[227,4] = "LP1912"
$ws1data[228,0] = "14:42:53"
$ws1data[228,1] = "15:10"
$ws1data[228,2] = "10_OLMOS"
$ws1data[228,3] = 28
$ws1data[228,4] = "LP1912"
$ws1data[229,0] = "14:57:45"
$ws1data[229,1] = "15:11"
$ws1data[229,2] = "10_OLMOS"
$ws1data[229,3] = 14
$ws1data[229,4] = "LP1912"
$ws1data[230,0] = "14:57:45"
$ws1data[230,1] = "15:18"
$ws1data[230,2] = "215C_EL PATO"
$ws1data[230,3] = 21
$ws1data[230,4] = "LP1912"
$ws1data[231,0] = "14:57:45"
$ws1data[231,1] = "15:23"
$ws1data[231,2] = "11_ETCHEVERRY"
$ws1data[231,3] = 26
$ws1data[231,4] = "LP1912"
$ws1data[232,0] = "14:57:45"
$ws1data[232,1] = "15:33"
$ws1data[232,2] = "23_HERNANDEZ"
$ws1data[232,3] = 36
$ws1data[232,4] = "LP1912"
$ws1data[233,0] = "14:57:45"
$ws1data[233,1] = "15:35"
$ws1data[233,2] = "10_OLMOS"
$ws1data[233,3] = 38
$ws1data[233,4] = "LP1912"
$ws1data[234,0] = "14:15:51"
$ws1data[234,1] = "15:36"
$ws1data[234,2] = "23_HERNANDEZ"
$ws1data[234,3] = 81
$ws1data[234,4] = "LP1912"
$ws1data[235,0] = "14:57:45"
$ws1data[235,1] = "15:38"
$ws1data[235,2] = "17X38_ROMERO"
$ws1data[235,3] = 41
$ws1data[235,4] = "LP1912"
$ws1data[236,0] = "14:42:53"
$ws1data[236,1] = "15:38"
$ws1data[236,2] = "23_HERNANDEZ"
$ws1data[236,3] = 56
$ws1data[236,4] = "LP1912"
$ws1data[237,0] = "14:57:45"
$ws1data[237,1] = "15:41"
$ws1data[237,2] = "16_SANTA ANA"
$ws1data[237,3] = 44
$ws1data[237,4] = "LP1912"
$ws1data[238,0] = "14:57:45"
$ws1data[238,1] = "15:47"
$ws1data[238,2] = "11_ETCHEVERRY"
$ws1data[238,3] = 50
$ws1data[238,4] = "LP1912"
$ws1data[239,0] = "14:57:45"
$ws1data[239,1] = "15:51"
$ws1data[239,2] = "16_P MOR-SANTA ANA"
$ws1data[239,3] = 54
$ws1data[239,4] = "LP1912"
$ws1data[240,0] = "14:42:53"
$ws1data[240,1] = "15:54"
$ws1data[240,2] = "27_EL RETIRO"
$ws1data[240,3] = 72
$ws1data[240,4] = "LP1912"
$ws1data[241,0] = "14:57:45"
$ws1data[241,1] = "15:54"
$ws1data[241,2] = "15_ABASTO"
$ws1data[241,3] = 57
$ws1data[241,4] = "LP1912"
$ws1data[242,0] = "14:57:45"
$ws1data[242,1] = "15:56"
$ws1data[242,2] = "27_EL RETIRO"
$ws1data[242,3] = 59
$ws1data[242,4] = "LP1912"
$ws1data[243,0] = "14:57:45"
$ws1data[243,1] = "15:59"
$ws1data[243,2] = "10_OLMOS"
$ws1data[243,3] = 62
$ws1data[243,4] = "LP1912"
$ws1data[244,0] = "14:42:53"
$ws1data[244,1] = "16:09"
$ws1data[244,2] = "15_ABASTO"
$ws1data[244,3] = 87
$ws1data[244,4] = "LP1912"
$ws1data[245,0] = "14:57:45"
$ws1data[245,1] = "16:10"
$ws1data[245,2] = "15_ABASTO"
$ws1data[245,3] = 73
$ws1data[245,4] = "LP1912"
$ws1data[246,0] = "14:57:45"
$ws1data[246,1] = "16:11"
$ws1data[246,2] = "14_ABASTO"
$ws1data[246,3] = 74
$ws1data[246,4] = "LP1912"
$ws1data[247,0] = "14:57:45"
$ws1data[247,1] = "16:23"
$ws1data[247,2] = "215_ALUAR"
$ws1data[247,3] = 86
$ws1data[247,4] = "LP1912"
$ws1data[248,0] = "14:57:45"
$ws1data[248,1] = "16:35"
$ws1data[248,2] = "11_ETCHEVERRY"
$ws1data[248,3] = 98
$ws1data[248,4] = "LP1912"
$ws1data[249,0] = "14:57:45"
$ws1data[249,1] = "16:35"
$ws1data[249,2] = "17X38_ROMERO"
$ws1data[249,3] = 98
$ws1data[249,4] = "LP1912"
$ws1data[250,0] = "14:57:45"
$ws1data[250,1] = "16:47"
$ws1data[250,2] = "225_GOMEZ"
$ws1data[250,3] = 110
$ws1data[250,4] = "LP1912"
$ws1data[251,0] = "14:57:45"
$ws1data[251,1] = "16:51"
$ws1data[251,2] = "16_P MOR-SANTA ANA"
$ws1data[251,3] = 114
$ws1data[251,4] = "LP1912"
$ws1.Range("A6:E257").Value2 = $ws1data

# --- LP1912-215 ---
$ws2.Cells.Item(2,1).Value = "Última actualización: 14:57:45"
$ws2.Cells.Item(3,1).Value = "Total filas: 43"
$ws2data = New-Object 'object[,]' 43,5
$ws2data[0,0] = "04:43:39"
$ws2data[0,1] = "04:45"
$ws2data[0,2] = "215A_EL PATO"
$ws2data[0,3] = 2
$ws2data[0,4] = "LP1912"
$ws2data[1,0] = "04:56:49"
$ws2data[1,1] = "05:34"
$ws2data[1,2] = "215B_EL PATO"
$ws2data[1,3] = 38
$ws2data[1,4] = "LP1912"
$ws2data[2,0] = "05:51:38"
$ws2data[2,1] = "06:10"
$ws2data[2,2] = "215A_EL PATO"
$ws2data[2,3] = 19
$ws2data[2,4] = "LP1912"
$ws2data[3,0] = "05:23:04"
$ws2data[3,1] = "06:11"
$ws2data[3,2] = "215A_EL PATO"
$ws2data[3,3] = 48
$ws2data[3,4] = "LP1912"
$ws2data[4,0] = "04:56:49"
$ws2data[4,1] = "06:46"
$ws2data[4,2] = "215C_EL PATO"
$ws2data[4,3] = 110
$ws2data[4,4] = "LP1912"
$ws2data[5,0] = "05:51:38"
$ws2data[5,1] = "06:50"
$ws2data[5,2] = "215A_EL PATO"
$ws2data[5,3] = 59
$ws2data[5,4] = "LP1912"
$ws2data[6,0] = "06:46:06"
$ws2data[6,1] = "06:51"
$ws2data[6,2] = "215A_EL PATO"
$ws2data[6,3] = 5
$ws2data[6,4] = "LP1912"
$ws2data[7,0] = "06:19:29"
$ws2data[7,1] = "07:06"
$ws2data[7,2] = "215C_EL PATO"
$ws2data[7,3] = 47
$ws2data[7,4] = "LP1912"
$ws2data[8,0] = "06:58:01"
$ws2data[8,1] = "07:07"
$ws2data[8,2] = "215C_EL PATO"
$ws2data[8,3] = 9
$ws2data[8,4] = "LP1912"
$ws2data[9,0] = "07:21:42"
$ws2data[9,1] = "07:21"
$ws2data[9,2] = "215A_EL PATO"
$ws2data[9,3] = 0
$ws2data[9,4] = "LP1912"
$ws2data[10,0] = "07:48:14"
$ws2data[10,1] = "07:48"
$ws2data[10,2] = "215A_EL PATO"
$ws2data[10,3] = 0
$ws2data[10,4] = "LP1912"
$ws2data[11,0] = "08:00:50"
$ws2data[11,1] = "08:00"
$ws2data[11,2] = "215A_EL PATO"
$ws2data[11,3] = 0
$ws2data[11,4] = "LP1912"
$ws2data[12,0] = "08:30:59"
$ws2data[12,1] = "08:30"
$ws2data[12,2] = "215A_EL PATO"
$ws2data[12,3] = 0
$ws2data[12,4] = "LP1912"
$ws2data[13,0] = "08:30:59"
$ws2data[13,1] = "08:33"
$ws2data[13,2] = "215C_EL PATO"
$ws2data[13,3] = 3
$ws2data[13,4] = "LP1912"
$ws2data[14,0] = "07:48:14"
$ws2data[14,1] = "08:34"
$ws2data[14,2] = "215C_EL PATO"
$ws2data[14,3] = 46
$ws2data[14,4] = "LP1912"
$ws2data[15,0] = "08:30:59"
$ws2data[15,1] = "08:47"
$ws2data[15,2] = "215A_EL PATO"
$ws2data[15,3] = 17
$ws2data[15,4] = "LP1912"
$ws2data[16,0] = "08:48:29"
$ws2data[16,1] = "08:48"
$ws2data[16,2] = "215A_EL PATO"
$ws2data[16,3] = 0
$ws2data[16,4] = "LP1912"
$ws2data[17,0] = "08:56:14"
$ws2data[17,1] = "08:56"
$ws2data[17,2] = "215A_EL PATO"
$ws2data[17,3] = 0
$ws2data[17,4] = "LP1912"
$ws2data[18,0] = "08:56:14"
$ws2data[18,1] = "08:56"
$ws2data[18,2] = "215B_EL PATO"
$ws2data[18,3] = 0
$ws2data[18,4] = "LP1912"
$ws2data[19,0] = "08:48:29"
$ws2data[19,1] = "08:59"
$ws2data[19,2] = "215B_EL PATO"
$ws2data[19,3] = 11
$ws2data[19,4] = "LP1912"
$ws2data[20,0] = "07:48:14"
$ws2data[20,1] = "09:00"
$ws2data[20,2] = "215B_EL PATO"
$ws2data[20,3] = 72
$ws2data[20,4] = "LP1912"
$ws2data[21,0] = "08:30:59"
$ws2data[21,1] = "09:18"
$ws2data[21,2] = "215_EL PELIGRO"
$ws2data[21,3] = 48
$ws2data[21,4] = "LP1912"
$ws2data[22,0] = "08:56:14"
$ws2data[22,1] = "09:19"
$ws2data[22,2] = "215_EL PELIGRO"
$ws2data[22,3] = 23
$ws2data[22,4] = "LP1912"
$ws2data[23,0] = "08:30:59"
$ws2data[23,1] = "10:03"
$ws2data[23,2] = "215C_EL PATO"
$ws2data[23,3] = 93
$ws2data[23,4] = "LP1912"
$ws2data[24,0] = "09:31:15"
$ws2data[24,1] = "10:04"
$ws2data[24,2] = "215C_EL PATO"
$ws2data[24,3] = 33
$ws2data[24,4] = "LP1912"
$ws2data[25,0] = "10:55:35"
$ws2data[25,1] = "11:19"
$ws2data[25,2] = "215C_EL PATO"
$ws2data[25,3] = 24
$ws2data[25,4] = "LP1912"
$ws2data[26,0] = "11:17:16"
$ws2data[26,1] = "11:33"
$ws2data[26,2] = "215A_EL PATO"
$ws2data[26,3] = 16
$ws2data[26,4] = "LP1912"
$ws2data[27,0] = "11:44:58"
$ws2data[27,1] = "11:44"
$ws2data[27,2] = "215B_EL PATO"
$ws2data[27,3] = 0
$ws2data[27,4] = "LP1912"
$ws2data[28,0] = "11:17:16"
$ws2data[28,1] = "11:45"
$ws2data[28,2] = "215B_EL PATO"
$ws2data[28,3] = 28
$ws2data[28,4] = "LP1912"
$ws2data[29,0] = "12:23:39"
$ws2data[29,1] = "12:33"
$ws2data[29,2] = "215C_EL PATO"
$ws2data[29,3] = 10
$ws2data[29,4] = "LP1912"
$ws2data[30,0] = "11:17:16"
$ws2data[30,1] = "12:34"
$ws2data[30,2] = "215C_EL PATO"
$ws2data[30,3] = 77
$ws2data[30,4] = "LP1912"
$ws2data[31,0] = "12:43:41"
$ws2data[31,1] = "13:03"
$ws2data[31,2] = "215C_EL PATO"
$ws2data[31,3] = 20
$ws2data[31,4] = "LP1912"
$ws2data[32,0] = "12:57:23"
$ws2data[32,1] = "13:04"
$ws2data[32,2] = "215C_EL PATO"
$ws2data[32,3] = 7
$ws2data[32,4] = "LP1912"
$ws2data[33,0] = "12:43:41"
$ws2data[33,1] = "13:18"
$ws2data[33,2] = "215_ALUAR"
$ws2data[33,3] = 35
$ws2data[33,4] = "LP1912"
$ws2data[34,0] = "12:57:23"
$ws2data[34,1] = "13:19"
$ws2data[34,2] = "215_ALUAR"
$ws2data[34,3] = 22
$ws2data[34,4] = "LP1912"
$ws2data[35,0] = "11:58:47"
$ws2data[35,1] = "13:37"
$ws2data[35,2] = "215A_EL PATO"
$ws2data[35,3] = 99
$ws2data[35,4] = "LP1912"
$ws2data[36,0] = "13:22:41"
$ws2data[36,1] = "13:38"
$ws2data[36,2] = "215A_EL PATO"
$ws2data[36,3] = 16
$ws2data[36,4] = "LP1912"
$ws2data[37,0] = "14:15:51"
$ws2data[37,1] = "14:18"
$ws2data[37,2] = "215C_EL PATO"
$ws2data[37,3] = 3
$ws2data[37,4] = "LP1912"
$ws2data[38,0] = "12:57:23"
$ws2data[38,1] = "14:19"
$ws2data[38,2] = "215C_EL PATO"
$ws2data[38,3] = 82
$ws2data[38,4] = "LP1912"
$ws2data[39,0] = "14:57:45"
$ws2data[39,1] = "14:58"
$ws2data[39,2] = "215B_EL PATO"
$ws2data[39,3] = 1
$ws2data[39,4] = "LP1912"
$ws2data[40,0] = "14:42:53"
$ws2data[40,1] = "14:59"
$ws2data[40,2] = "215B_EL PATO"
$ws2data[40,3] = 17
$ws2data[40,4] = "LP1912"
$ws2data[41,0] = "14:57:45"
$ws2data[41,1] = "15:18"
$ws2data[41,2] = "215C_EL PATO"
$ws2data[41,3] = 21
$ws2data[41,4] = "LP1912"
$ws2data[42,0] = "14:57:45"
$ws2data[42,1] = "16:23"
$ws2data[42,2] = "215_ALUAR"
$ws2data[42,3] = 86
$ws2data[42,4] = "LP1912"
$ws2.Range("A6:E48").Value2 = $ws2data

# --- 6203-6173 ---
$ws3.Cells.Item(2,1).Value = "Última actualización: 14:57:45"
$ws3.Cells.Item(3,1).Value = "Total filas: 45"
$ws3data = New-Object 'object[,]' 45,5
$ws3data[0,0] = "04:56:49"
$ws3data[0,1] = "05:43"
$ws3data[0,2] = "215A_LA PLATA"
$ws3data[0,3] = 47
$ws3data[0,4] = "L6173"
$ws3data[1,0] = "04:56:49"
$ws3data[1,1] = "06:08"
$ws3data[1,2] = "215A_LA PLATA"
$ws3data[1,3] = 72
$ws3data[1,4] = "L6173"
$ws3data[2,0] = "04:56:49"
$ws3data[2,1] = "06:32"
$ws3data[2,2] = "215C_LA PLATA"
$ws3data[2,3] = 96
$ws3data[2,4] = "L6203"
$ws3data[3,0] = "06:19:29"
$ws3data[3,1] = "07:27"
$ws3data[3,2] = "215A_LA PLATA"
$ws3data[3,3] = 68
$ws3data[3,4] = "L6173"
$ws3data[4,0] = "06:46:06"
$ws3data[4,1] = "07:32"
$ws3data[4,2] = "215A_LA PLATA"
$ws3data[4,3] = 46
$ws3data[4,4] = "L6173"
$ws3data[5,0] = "07:21:42"
$ws3data[5,1] = "07:35"
$ws3data[5,2] = "215A_LA PLATA"
$ws3data[5,3] = 14
$ws3data[5,4] = "L6173"
$ws3data[6,0] = "06:58:01"
$ws3data[6,1] = "07:36"
$ws3data[6,2] = "215A_LA PLATA"
$ws3data[6,3] = 38
$ws3data[6,4] = "L6173"
$ws3data[7,0] = "08:00:50"
$ws3data[7,1] = "08:09"
$ws3data[7,2] = "215A_LA PLATA"
$ws3data[7,3] = 9
$ws3data[7,4] = "L6173"
$ws3data[8,0] = "07:48:14"
$ws3data[8,1] = "08:10"
$ws3data[8,2] = "215A_LA PLATA"
$ws3data[8,3] = 22
$ws3data[8,4] = "L6173"
$ws3data[9,0] = "07:21:42"
$ws3data[9,1] = "08:23"
$ws3data[9,2] = "215C_LA PLATA"
$ws3data[9,3] = 62
$ws3data[9,4] = "L6203"
$ws3data[10,0] = "08:00:50"
$ws3data[10,1] = "08:26"
$ws3data[10,2] = "215C_LA PLATA"
$ws3data[10,3] = 26
$ws3data[10,4] = "L6203"
$ws3data[11,0] = "06:58:01"
$ws3data[11,1] = "08:52"
$ws3data[11,2] = "215A_LA PLATA"
$ws3data[11,3] = 114
$ws3data[11,4] = "L6173"
$ws3data[12,0] = "07:48:14"
$ws3data[12,1] = "09:09"
$ws3data[12,2] = "215A_LA PLATA"
$ws3data[12,3] = 81
$ws3data[12,4] = "L6173"
$ws3data[13,0] = "08:00:50"
$ws3data[13,1] = "09:20"
$ws3data[13,2] = "215A_LA PLATA"
$ws3data[13,3] = 80
$ws3data[13,4] = "L6173"
$ws3data[14,0] = "08:30:59"
$ws3data[14,1] = "09:50"
$ws3data[14,2] = "215A_LA PLATA"
$ws3data[14,3] = 80
$ws3data[14,4] = "L6173"
$ws3data[15,0] = "08:30:59"
$ws3data[15,1] = "09:55"
$ws3data[15,2] = "215C_LA PLATA"
$ws3data[15,3] = 85
$ws3data[15,4] = "L6203"
$ws3data[16,0] = "09:31:15"
$ws3data[16,1] = "09:56"
$ws3data[16,2] = "215C_LA PLATA"
$ws3data[16,3] = 25
$ws3data[16,4] = "L6203"
$ws3data[17,0] = "08:48:29"
$ws3data[17,1] = "10:08"
$ws3data[17,2] = "215A_LA PLATA"
$ws3data[17,3] = 80
$ws3data[17,4] = "L6173"
$ws3data[18,0] = "08:30:59"
$ws3data[18,1] = "10:10"
$ws3data[18,2] = "215A_LA PLATA"
$ws3data[18,3] = 100
$ws3data[18,4] = "L6173"
$ws3data[19,0] = "08:56:14"
$ws3data[19,1] = "10:11"
$ws3data[19,2] = "215A_LA PLATA"
$ws3data[19,3] = 75
$ws3data[19,4] = "L6173"
$ws3data[20,0] = "08:56:14"
$ws3data[20,1] = "10:16"
$ws3data[20,2] = "215A_LA PLATA"
$ws3data[20,3] = 80
$ws3data[20,4] = "L6173"
$ws3data[21,0] = "09:31:15"
$ws3data[21,1] = "10:18"
$ws3data[21,2] = "215A_LA PLATA"
$ws3data[21,3] = 47
$ws3data[21,4] = "L6173"
$ws3data[22,0] = "08:30:59"
$ws3data[22,1] = "10:21"
$ws3data[22,2] = "215B_LP-P MOR-1 Y 57"
$ws3data[22,3] = 111
$ws3data[22,4] = "L6173"
$ws3data[23,0] = "09:31:15"
$ws3data[23,1] = "10:22"
$ws3data[23,2] = "215B_LP-P MOR-1 Y 57"
$ws3data[23,3] = 51
$ws3data[23,4] = "L6173"
$ws3data[24,0] = "11:44:58"
$ws3data[24,1] = "11:55"
$ws3data[24,2] = "215C_LA PLATA"
$ws3data[24,3] = 11
$ws3data[24,4] = "L6203"
$ws3data[25,0] = "11:17:16"
$ws3data[25,1] = "11:56"
$ws3data[25,2] = "215C_LA PLATA"
$ws3data[25,3] = 39
$ws3data[25,4] = "L6203"
$ws3data[26,0] = "12:43:41"
$ws3data[26,1] = "12:55"
$ws3data[26,2] = "215C_LA PLATA"
$ws3data[26,3] = 12
$ws3data[26,4] = "L6203"
$ws3data[27,0] = "11:17:16"
$ws3data[27,1] = "12:56"
$ws3data[27,2] = "215C_LA PLATA"
$ws3data[27,3] = 99
$ws3data[27,4] = "L6203"
$ws3data[28,0] = "12:43:41"
$ws3data[28,1] = "13:10"
$ws3data[28,2] = "215A_LA PLATA"
$ws3data[28,3] = 27
$ws3data[28,4] = "L6173"
$ws3data[29,0] = "12:57:23"
$ws3data[29,1] = "13:11"
$ws3data[29,2] = "215A_LA PLATA"
$ws3data[29,3] = 14
$ws3data[29,4] = "L6173"
$ws3data[30,0] = "12:43:41"
$ws3data[30,1] = "13:21"
$ws3data[30,2] = "215B_LP-P MOR-1 Y 57"
$ws3data[30,3] = 38
$ws3data[30,4] = "L6173"
$ws3data[31,0] = "13:22:41"
$ws3data[31,1] = "13:22"
$ws3data[31,2] = "215B_LP-P MOR-1 Y 57"
$ws3data[31,3] = 0
$ws3data[31,4] = "L6173"
$ws3data[32,0] = "12:43:41"
$ws3data[32,1] = "13:55"
$ws3data[32,2] = "215C_LA PLATA"
$ws3data[32,3] = 72
$ws3data[32,4] = "L6203"
$ws3data[33,0] = "12:57:23"
$ws3data[33,1] = "13:56"
$ws3data[33,2] = "215C_LA PLATA"
$ws3data[33,3] = 59
$ws3data[33,4] = "L6203"
$ws3data[34,0] = "13:51:40"
$ws3data[34,1] = "13:57"
$ws3data[34,2] = "215C_LA PLATA"
$ws3data[34,3] = 6
$ws3data[34,4] = "L6203"
$ws3data[35,0] = "13:22:41"
$ws3data[35,1] = "14:00"
$ws3data[35,2] = "215C_LA PLATA"
$ws3data[35,3] = 38
$ws3data[35,4] = "L6203"
$ws3data[36,0] = "14:42:53"
$ws3data[36,1] = "14:45"
$ws3data[36,2] = "215C_LA PLATA"
$ws3data[36,3] = 3
$ws3data[36,4] = "L6203"
$ws3data[37,0] = "12:57:23"
$ws3data[37,1] = "14:46"
$ws3data[37,2] = "215C_LA PLATA"
$ws3data[37,3] = 109
$ws3data[37,4] = "L6203"
$ws3data[38,0] = "13:51:40"
$ws3data[38,1] = "15:05"
$ws3data[38,2] = "215A_LA PLATA"
$ws3data[38,3] = 74
$ws3data[38,4] = "L6173"
$ws3data[39,0] = "14:15:51"
$ws3data[39,1] = "15:06"
$ws3data[39,2] = "215A_LA PLATA"
$ws3data[39,3] = 51
$ws3data[39,4] = "L6173"
$ws3data[40,0] = "14:42:53"
$ws3data[40,1] = "15:09"
$ws3data[40,2] = "215A_LA PLATA"
$ws3data[40,3] = 27
$ws3data[40,4] = "L6173"
$ws3data[41,0] = "14:57:45"
$ws3data[41,1] = "15:11"
$ws3data[41,2] = "215A_LA PLATA"
$ws3data[41,3] = 14
$ws3data[41,4] = "L6173"
$ws3data[42,0] = "14:57:45"
$ws3data[42,1] = "15:45"
$ws3data[42,2] = "215C_LA PLATA"
$ws3data[42,3] = 48
$ws3data[42,4] = "L6203"
$ws3data[43,0] = "14:57:45"
$ws3data[43,1] = "16:42"
$ws3data[43,2] = "215B_LP-P MOR-1 Y 57"
$ws3data[43,3] = 105
$ws3data[43,4] = "L6173"
$ws3data[44,0] = "14:57:45"
$ws3data[44,1] = "16:55"
$ws3data[44,2] = "215C_LA PLATA"
$ws3data[44,3] = 118
$ws3data[44,4] = "L6203"
$ws3.Range("A6:E50").Value2 = $ws3data

Write-Output "Updated LP1912 ($($ws1.UsedRange.Rows.Count) rows), LP1912-215, 6203-6173"
